# Apply needle calibration data correction: the data rows (2-8) need to be
# re-sorted in ascending order of column A (time), which results in the
# following reordering (1-indexed data rows as they appeared originally):
#   2 <- 5, 3 <- 2, 4 <- 4, 5 <- 3, 6 <- 6, 7 <- 8, 8 <- 7

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow = 8
$numCols = 4

# Snapshot the current values of the data rows before any writes, so that
# later writes don't clobber values we still need to read.
$original = @{}
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $rowVals = @()
    for ($c = 1; $c -le $numCols; $c++) {
        $rowVals += , ($ws.Cells.Item($r, $c).Value2)
    }
    $original[$r] = $rowVals
}

# Mapping of destination row -> source row (time-sorted order).
$rowMap = @{
    2 = 5
    3 = 2
    4 = 4
    5 = 3
    6 = 6
    7 = 8
    8 = 7
}

foreach ($destRow in ($rowMap.Keys | Sort-Object)) {
    $srcRow = $rowMap[$destRow]
    $srcVals = $original[$srcRow]
    for ($c = 1; $c -le $numCols; $c++) {
        $ws.Cells.Item($destRow, $c).Value = $srcVals[$c - 1]
    }
}
